# Add DCL revoke and reset password case:
# - Change B2 and B3 from "n" to "y" (Testable column)
# - Move selection to C11

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "y"
$ws.Range("B3").Value = "y"

$ws.Range("C11").Select()
